# Add GBR3 rows to the derived data map (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 17: GBR3 / ageband
$ws.Cells.Item(17, 1).Value = "GBR3"
$ws.Cells.Item(17, 2).Value = "ageband"
$ws.Cells.Item(17, 3).Value = "data/derived/GBR3/GBR3_agebands.RDS"
$ws.Cells.Item(17, 4).Value = "marginal"
$ws.Cells.Item(17, 5).Value = "aggregate"

# New row 18: GBR3 / region
$ws.Cells.Item(18, 1).Value = "GBR3"
$ws.Cells.Item(18, 2).Value = "region"
$ws.Cells.Item(18, 3).Value = "data/derived/GBR3/GBR3_regions.RDS"
$ws.Cells.Item(18, 4).Value = "marginal"
$ws.Cells.Item(18, 5).Value = "aggregate"

# Match the author's final cell selection after entering the new data
$ws.Range("C19").Select()
